{"js": "// HREOS-2824: Fwd: RV: Nuevo circuito reservas Sareb.pptx\n//\n// The \"COBRO RESERVA\" paragraph that lists the payment-reference\n// parameters is updated:\n//   - \"y (iii)\" -> \"(iii)\"                      (drop the stray \"y\")\n//   - \"(iii) d.n.i.\" -> \"(iii) n\u00famero de oferta; (iv)D.N.I\"\n//        (a new \"n\u00famero de oferta\" item (iii) is inserted, and the old\n//         \"d.n.i.\" item is renumbered to (iv) and re-cased to \"D.N.I\")\n// The cursor (the \"_GoBack\" last-edit-position bookmark) ends up right\n// before the untouched trailing \"de los compradores.\" text.\n\nconst body = context.document.body;\n\n// 1) Drop the \"y \" that used to precede \"(iii)\".\nconst yMatches = body.search(\"y (iii) \", { matchCase: false, ignorePunct: false });\nyMatches.load(\"text\");\nawait context.sync();\n\nif (yMatches.items.length !== 1) {\n  throw new Error(\n    `Expected exactly one \"y (iii) \" match, found ${yMatches.items.length}`\n  );\n}\nyMatches.items[0].insertText(\"(iii) \", \"Replace\");\nawait context.sync();\n\n// 2) Turn the old \"(iii) d.n.i.\" item into the new \"(iii) n\u00famero de\n//    oferta; (iv)D.N.I\" pair of items.\nconst dniMatches = body.search(\"d.n.i.\", { matchCase: false, ignorePunct: false });\ndniMatches.load(\"text\");\nawait context.sync();\n\nif (dniMatches.items.length !== 1) {\n  throw new Error(\n    `Expected exactly one \"d.n.i.\" match, found ${dniMatches.items.length}`\n  );\n}\ndniMatches.items[0].insertText(\"n\u00famero de oferta; (iv)D.N.I\", \"Replace\");\nawait context.sync();\n\n// 3) Re-anchor the \"_GoBack\" bookmark (Word's last-edit-position marker)\n//    to sit right before the untouched \"de los compradores.\" tail, which\n//    is where the author's cursor was left after typing the replacement\n//    text above.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst tailMatches = body.search(\"de los compradores. \", {\n  matchCase: false,\n  ignorePunct: false,\n});\ntailMatches.load(\"text\");\nawait context.sync();\n\nif (tailMatches.items.length !== 1) {\n  throw new Error(\n    `Expected exactly one \"de los compradores. \" match, found ${tailMatches.items.length}`\n  );\n}\ntailMatches.items[0].getRange(\"Start\").insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# HREOS-2824: Fwd: RV: Nuevo circuito reservas Sareb.pptx\n#\n# The \"COBRO RESERVA\" paragraph that lists the payment-reference\n# parameters is updated:\n#   - \"y (iii)\" -> \"(iii)\"                      (drop the stray \"y\")\n#   - \"(iii) d.n.i.\" -> \"(iii) n\u00famero de oferta; (iv)D.N.I\"\n#        (a new \"n\u00famero de oferta\" item (iii) is inserted, and the old\n#         \"d.n.i.\" item is renumbered to (iv) and re-cased to \"D.N.I\")\n# The cursor (the \"_GoBack\" last-edit-position bookmark) ends up right\n# before the untouched trailing \"de los compradores.\" text.\n\n$d = $word.ActiveDocument\n\n# 1) Drop the \"y \" that used to precede \"(iii)\".\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Replacement.ClearFormatting()\n$found1 = $rng1.Find.Execute(\"y (iii) \", $false, $false, $false, $false, $false, $true, 1, $false, \"(iii) \", 2)\nif (-not $found1) {\n    throw \"Could not find 'y (iii) ' to replace\"\n}\n\n# 2) Turn the old \"(iii) d.n.i.\" item into the new \"(iii) n\u00famero de\n#    oferta; (iv)D.N.I\" pair of items.\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Replacement.ClearFormatting()\n$found2 = $rng2.Find.Execute(\"d.n.i.\", $false, $false, $false, $false, $false, $true, 1, $false, \"n\u00famero de oferta; (iv)D.N.I\", 2)\nif (-not $found2) {\n    throw \"Could not find 'd.n.i.' to replace\"\n}\n\n# 3) Re-anchor the \"_GoBack\" bookmark (Word's last-edit-position marker)\n#    to sit right before the untouched \"de los compradores.\" tail, which\n#    is where the author's cursor was left after typing the replacement\n#    text above.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n$rng3 = $d.Content\n$found3 = $rng3.Find.Execute(\"de los compradores. \")\nif (-not $found3) {\n    throw \"Could not find 'de los compradores. ' tail\"\n}\n$rng3.Collapse(1)\n$d.Bookmarks.Add(\"_GoBack\", $rng3)\n"}
